$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 0.8199996238250412
$ws.Range("C4").Value = 0.04243367311486969
$ws.Range("I4").Value = 0.6929483976074302
$ws.Range("J4").Value = 0.5990338164251208
$ws.Range("O4").Value = 0.750602870813397
$ws.Range("P4").Value = 0.07018039593915391
$ws.Range("B5").Value = 0.780408163265306
$ws.Range("C5").Value = 0.009795918367346968
$ws.Range("D5").Value = 0.7012987012987013
$ws.Range("E5").Value = 0.8544303797468354
$ws.Range("F5").Value = 0.54
$ws.Range("G5").Value = 0.54
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0.7127577336849282
$ws.Range("J5").Value = 0.6328502415458936
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 23
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 27
$ws.Range("O5").Value = 0.6768421052631579
$ws.Range("P5").Value = 0.05960854575975858
$ws.Range("B6").Value = 0.7291809972169407
$ws.Range("C6").Value = 0.05573953308217269
$ws.Range("D6").Value = 0.5762711864406779
$ws.Range("E6").Value = 0.6071428571428571
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.53125
$ws.Range("H6").Value = 0.6296296296296297
$ws.Range("I6").Value = 0.6101047952542277
$ws.Range("J6").Value = 0.5652173913043479
$ws.Range("K6").Value = 8
$ws.Range("L6").Value = 15
$ws.Range("M6").Value = 10
$ws.Range("N6").Value = 17
$ws.Range("O6").Value = 0.6249090909090909
$ws.Range("P6").Value = 0.04878340505235926
$ws.Range("B7").Value = 0.7230060072562898
$ws.Range("C7").Value = 0.07565919763686081
$ws.Range("D7").Value = 0.6885245901639345
$ws.Range("E7").Value = 0.73943661971831
$ws.Range("F7").Value = 0.62
$ws.Range("G7").Value = 0.6176470588235294
$ws.Range("H7").Value = 0.7777777777777778
$ws.Range("I7").Value = 0.5850495080186957
$ws.Range("J7").Value = 0.5869565217391305
$ws.Range("K7").Value = 10
$ws.Range("L7").Value = 13
$ws.Range("M7").Value = 6
$ws.Range("N7").Value = 21
$ws.Range("O7").Value = 0.5402775119617225
$ws.Range("P7").Value = 0.0787238965615758
$ws.Range("B8").Value = 0.7801252236135958
$ws.Range("C8").Value = 0.03863755124112093
$ws.Range("D8").Value = 0.5333333333333333
$ws.Range("E8").Value = 0.5673758865248226
$ws.Range("F8").Value = 0.44
$ws.Range("G8").Value = 0.4848484848484849
$ws.Range("H8").Value = 0.5925925925925926
$ws.Range("I8").Value = 0.5877398214627141
$ws.Range("J8").Value = 0.4718196457326893
$ws.Range("K8").Value = 6
$ws.Range("L8").Value = 17
$ws.Range("M8").Value = 11
$ws.Range("N8").Value = 16
$ws.Range("O8").Value = 0.7413875598086124
$ws.Range("P8").Value = 0.04884368388503456
